# Regenerate Handback status report: refresh the "Latest HO Xliff Generate
# Date" / "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the a6f5cf38-7a53-43f0-9081-cb5230a1fab6 file (row 3 on
# each sheet) to reflect the newly generated handback report.

$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date (column G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-22 16:49:22"

# zh-cn sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-22 16:49:17"
$wsZhCn.Range("K3").Value = "2016-08-22 16:49:40"

# de-de sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-22 16:49:22"
$wsDeDe.Range("K3").Value = "2016-08-22 16:49:49"
